# Applies the edits described in the commit diff to the single worksheet.
# Every touched cell holds a literal cached value (no formulas in this
# workbook), so each row-level quantity/value change, row swap, subtotal
# and grand-total cascade is written explicitly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23
$ws.Range("F23").Value = 39
$ws.Range("G23").Value = 1598.61

# Row 34
$ws.Range("B34").Value = 68990.94

# Row 38
$ws.Range("F38").Value = 533
$ws.Range("G38").Value = 19406.53

# Row 42
$ws.Range("F42").Value = 85
$ws.Range("G42").Value = 3580.2

# Row 45
$ws.Range("F45").Value = 94
$ws.Range("G45").Value = 2137.56

# Row 53
$ws.Range("F53").Value = 43
$ws.Range("G53").Value = 705.63

# Row 58
$ws.Range("F58").Value = 99
$ws.Range("G58").Value = 7715.07

# Row 66
$ws.Range("B66").Value = 250280.94

# Row 95
$ws.Range("F95").Value = 7
$ws.Range("G95").Value = 2602.6

# Row 97
$ws.Range("B97").Value = 20809.69

# Row 100
$ws.Range("F100").Value = 17
$ws.Range("G100").Value = 2924.68

# Row 123
$ws.Range("B123").Value = 86647.09

# Row 161
$ws.Range("B161").Value = 57756
$ws.Range("E161").Value = 79.37
$ws.Range("F161").Value = -100
$ws.Range("G161").Value = -6644

# Row 162
$ws.Range("B162").Value = 64350
$ws.Range("E162").Value = 70.63
$ws.Range("F162").Value = 2
$ws.Range("G162").Value = 132.88

# Row 212
$ws.Range("F212").Value = 96
$ws.Range("G212").Value = 8552.639999999999

# Row 213
$ws.Range("F213").Value = 236
$ws.Range("G213").Value = 29896.48

# Row 218
$ws.Range("B218").Value = 95279.35000000001

# Row 232
$ws.Range("F232").Value = 34
$ws.Range("G232").Value = 3897.42

# Row 240
$ws.Range("B240").Value = 16413.31

# Row 264
$ws.Range("F264").Value = 128
$ws.Range("G264").Value = 4459.52

# Row 278
$ws.Range("F278").Value = 59
$ws.Range("G278").Value = 7996.86

# Row 287
$ws.Range("F287").Value = 90
$ws.Range("G287").Value = 4926.6

# Row 290
$ws.Range("B290").Value = 64983
$ws.Range("C290").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F290").Value = 6
$ws.Range("G290").Value = 514.08

# Row 291
$ws.Range("B291").Value = 66194
$ws.Range("C291").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F291").Value = 35
$ws.Range("G291").Value = 2998.8

# Row 292
$ws.Range("B292").Value = 64985
$ws.Range("C292").Value = "HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S"
$ws.Range("F292").Value = 13
$ws.Range("G292").Value = 1140.1

# Row 293
$ws.Range("B293").Value = 66196
$ws.Range("C293").Value = "HIM-Total Care Baby Pants Drapers-Xl-9S"
$ws.Range("F293").Value = 27
$ws.Range("G293").Value = 2367.9

# Row 295
$ws.Range("B295").Value = 160106.17

# Row 308
$ws.Range("B308").Value = 63510
$ws.Range("E308").Value = 50.66
$ws.Range("F308").Value = 84
$ws.Range("G308").Value = 4001.76

# Row 309
$ws.Range("B309").Value = 55356
$ws.Range("E309").Value = 54.04
$ws.Range("F309").Value = -158
$ws.Range("G309").Value = -7527.12

# Row 344
$ws.Range("F344").Value = 12
$ws.Range("G344").Value = 2029.68

# Row 356
$ws.Range("B356").Value = 88358.69

# Row 371
$ws.Range("F371").Value = 85
$ws.Range("G371").Value = 12772.95

# Row 372
$ws.Range("B372").Value = 78130.87

# Row 381
$ws.Range("B381").Value = 47097
$ws.Range("D381").Value = 112.28
$ws.Range("E381").Value = 134.16
$ws.Range("F381").Value = 15
$ws.Range("G381").Value = 1684.2

# Row 382
$ws.Range("B382").Value = 58047
$ws.Range("D382").Value = 105.54
$ws.Range("E382").Value = 126.1
$ws.Range("F382").Value = 33
$ws.Range("G382").Value = 3482.82

# Row 387
$ws.Range("F387").Value = 516
$ws.Range("G387").Value = 49845.6

# Row 389
$ws.Range("B389").Value = 67719.91

# Row 408
$ws.Range("F408").Value = 261
$ws.Range("G408").Value = 4136.85

# Row 415
$ws.Range("F415").Value = 69
$ws.Range("G415").Value = 3774.3

# Row 417
$ws.Range("B417").Value = 192704.96

# Row 433
$ws.Range("F433").Value = 192
$ws.Range("G433").Value = 1850.88

# Row 438
$ws.Range("B438").Value = 32704.45

# Row 569
$ws.Range("F569").Value = 7
$ws.Range("G569").Value = 4093.04

# Row 573
$ws.Range("B573").Value = 42754

# Row 609
$ws.Range("F609").Value = 75
$ws.Range("G609").Value = 8160.75

# Row 611
$ws.Range("F611").Value = 9
$ws.Range("G611").Value = 766.8

# Row 623
$ws.Range("F623").Value = 92
$ws.Range("G623").Value = 47337.68

# Row 626
$ws.Range("F626").Value = 41
$ws.Range("G626").Value = 1935.61

# Row 628
$ws.Range("B628").Value = 256765.33

# Row 665
$ws.Range("F665").Value = 4
$ws.Range("G665").Value = 178.48

# Row 668
$ws.Range("B668").Value = 18522.98

# Row 673
$ws.Range("F673").Value = 53
$ws.Range("G673").Value = 1602.19

# Row 674
$ws.Range("F674").Value = 1373
$ws.Range("G674").Value = 223950.03

# Row 676
$ws.Range("F676").Value = 34
$ws.Range("G676").Value = 9617.58

# Row 677
$ws.Range("F677").Value = 44
$ws.Range("G677").Value = 6364.6

# Row 680
$ws.Range("B680").Value = 242291.76

# Row 696
$ws.Range("F696").Value = 10
$ws.Range("G696").Value = 4367

# Row 713
$ws.Range("B713").Value = 83840.92999999999

# Row 718
$ws.Range("B718").Value = 3661296.75

# Row 719
$ws.Range("B719").Value = 3661296.75
